# Append one new data row (row 88) to the daily log sheet, matching the
# existing rows' layout: date (text), weekday (text), hour (number),
# ranking (number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to be treated as text so the date-like string
# "2025/10/10" is stored literally (as the other date cells in the
# sheet are) instead of being auto-converted into an Excel date serial
# number. Resetting the style back to "Normal" afterwards keeps the
# cell format identical to its neighbours (no explicit style index).
$ws.Range("A88").NumberFormat = "@"
$ws.Range("A88").Value = "2025/10/10"
$ws.Range("A88").Style = "Normal"

$ws.Range("B88").Value = "金"
$ws.Range("C88").Value = 6
$ws.Range("D88").Value = 201
